# Apply "Added in gpt-5 results" edit to ResultFigures.xlsx
#
# Net effect of the change: the J column of the "EasternDataBQA" block
# (rows 6, 12, 18 and 24) currently reads "GPT-5-mini" and is being
# relabeled to "GPT-5(-mini)". The neighbouring "GPT-4.1-mini" rows
# (7, 13, 19, 25) keep their text as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ResultFigures")

# Rename the GPT-5-mini label to GPT-5(-mini) everywhere it appears
# in the EasternDataBQA results table (column J).
$ws.Range("J6").Value  = "GPT-5(-mini)"
$ws.Range("J12").Value = "GPT-5(-mini)"
$ws.Range("J18").Value = "GPT-5(-mini)"
$ws.Range("J24").Value = "GPT-5(-mini)"

# Match the recorded window/selection state: scroll back to the top of
# the sheet (default top-left cell) and leave the active selection on
# K23 instead of the previous K14.
$ws.Range("A1").Select() | Out-Null
$ws.Range("K23").Select() | Out-Null

# Reflect the resized/repositioned (maximized) application window
# recorded in the saved workbook view - best effort, harmless no-op if
# the host does not surface window geometry back into the OOXML.
$win = $wb.Windows.Item(1)
$win.WindowState = -4137
$win.Left   = -120
$win.Top    = -120
$win.Width  = 38640
$win.Height = 21240

$wb.Save()
